# fix typo + add option to skip prt linking
#
# This reproduces the diff that:
#  - adds 3 new shared strings (VTC.LINK_PRT / the recommendation note / "Link PRT to VTCs")
#  - fills in the previously-blank row 30 of the settings table with a new
#    "Link PRT to VTCs" option (default TRUE), wraps + taller-rows the Comments cell
#  - grows Table1 (and its AutoFilter) from A1:E29 to A1:E30
#  - moves the saved scroll position / selection down to the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the new row 30 --------------------------------------------
# Field_ID, then Comments, then Description - matching the order the
# strings were originally typed in (this controls shared-string order).
$ws.Range("E30").Value = "VTC.LINK_PRT"
$ws.Range("D30").Value = "This step is recommended, but you could disable it to save time if the linking has already been performed."
$ws.Range("B30").Value = "Link PRT to VTCs"
$ws.Range("C30").Value = $true

# Wrap the long Comments text and make the row tall enough to show it.
# (Touching Font explicitly mirrors the workbook's new dedicated cell style
# that was introduced for this cell - applyFont + wrapText, non-bold.)
$ws.Range("D30").WrapText = $true
$ws.Range("D30").Font.ThemeColor = 1
$ws.Rows.Item(30).RowHeight = 30

# --- Grow the table / autofilter to include the new row -----------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E30"))

# --- Update the saved view state (scroll position + selection) ----------
$ws.Application.GoTo($ws.Range("A13"), $true)
[void]$ws.Range("C30").Select()
